# "hall effect and capactive datasets recorded"
#
# Row 12 = Hall effect (code H1): the "Sudong?" tentative/unconfirmed
# assignment (shown in italics) is now confirmed, so it becomes a plain
# "Sudong", and the "Exists?" flag cell moves from the red (not-yet) fill
# to the green (exists) fill used throughout the rest of the sheet.
#
# Row 13 = Capacitive? (code C1): the "Qinghua" assignment had a stray
# leftover style; normalise it to the plain/default formatting used by
# every other row now that this dataset is recorded too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hall effect (row 12): confirm "Sudong?" -> "Sudong", drop the italics
$ws.Range("D12").Value = "Sudong"
$ws.Range("D12").Font.Italic = $false

# Hall effect now exists -> green fill (same colour already used on E3 etc.)
$ws.Range("E12").Interior.Color = $ws.Range("E3").Interior.Color

# Capacitive (row 13): clear the stray formatting on the People Involved cell
$ws.Range("D13").Font.Italic = $false

# Leave the cursor where the author ended up
$ws.Range("D22").Select()
